$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 148. This shifts the existing rows 148:170 down to
# 149:171 (carrying their formatting/content along, as Excel does natively).
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new weekly data point.
$ws.Cells.Item(148, 1).Value = 4
$ws.Cells.Item(148, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(148, 3).Value = "Los Lagos"
$ws.Cells.Item(148, 4).Value = 44505
$ws.Cells.Item(148, 5).Value = 10
$ws.Cells.Item(148, 6).Value = 100112043
$ws.Cells.Item(148, 7).Value = "Pepino ensalada"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 100
$ws.Cells.Item(148, 11).Value = 12000
$ws.Cells.Item(148, 12).Value = 12000
$ws.Cells.Item(148, 13).Value = 12000
$ws.Cells.Item(148, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(148, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(148, 16).Value = 200
$ws.Cells.Item(148, 17).Value = 60
$ws.Cells.Item(148, 18).Value = "Hortaliza"
